# =========================================================================
# Edit: add a new "2022-Q4" sheet (inserted right after the "总计" summary
# sheet) with fund-holding data, and update the "总计" summary table with
# a new leading row for 2022-Q4 (shifting the existing rows down by one).
# =========================================================================

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4
#    and renumber the index column (A) for the rows that shift down.
# -------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 20
$summary.Range("D2").Value = 7.01

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

# -------------------------------------------------------------------
# 2) Insert a new worksheet named "2022-Q4" right after "总计" (so it
#    lands before the existing "2022-Q3" sheet, pushing the rest down
#    one position, same as the target workbook).
# -------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Borrow the header-row / index-column formatting from the "2022-Q3"
# sheet so the new sheet's styling (bold header, bordered index column)
# matches the rest of the workbook.
$srcSheet = $wb.Worksheets.Item("2022-Q3")
$srcSheet.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$srcSheet.Range("A2").Copy()
$newSheet.Range("A2:A21").PasteSpecial(-4122)

# Force the numeric-looking text columns to be stored as text (matches
# the source data, which keeps trailing zeros like "102.10"/"24.17").
$newSheet.Range("B2:B21").NumberFormat = "@"
$newSheet.Range("D2:G21").NumberFormat = "@"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows
$newSheet.Range('B2').Value = '009100'
$newSheet.Range('C2').Value = '安信稳健增利混合A'
$newSheet.Range('D2').Value = '102.10'
$newSheet.Range('E2').Value = '24.17'
$newSheet.Range('F2').Value = '1.91'
$newSheet.Range('G2').Value = '1.9501'
$newSheet.Range('A2').Value = 0
$newSheet.Range('H2').Value = 7
$newSheet.Range('B3').Value = '009101'
$newSheet.Range('C3').Value = '安信稳健增利混合C'
$newSheet.Range('D3').Value = '67.40'
$newSheet.Range('E3').Value = '24.17'
$newSheet.Range('F3').Value = '1.91'
$newSheet.Range('G3').Value = '1.2873'
$newSheet.Range('A3').Value = 1
$newSheet.Range('H3').Value = 7
$newSheet.Range('B4').Value = '010864'
$newSheet.Range('C4').Value = '泓德卓远混合A'
$newSheet.Range('D4').Value = '25.76'
$newSheet.Range('E4').Value = '92.60'
$newSheet.Range('F4').Value = '2.83'
$newSheet.Range('G4').Value = '0.7290'
$newSheet.Range('A4').Value = 2
$newSheet.Range('H4').Value = 10
$newSheet.Range('B5').Value = '012256'
$newSheet.Range('C5').Value = '安信丰穗一年持有混合A'
$newSheet.Range('D5').Value = '24.73'
$newSheet.Range('E5').Value = '25.87'
$newSheet.Range('F5').Value = '1.81'
$newSheet.Range('G5').Value = '0.4476'
$newSheet.Range('A5').Value = 3
$newSheet.Range('H5').Value = 7
$newSheet.Range('B6').Value = '009849'
$newSheet.Range('C6').Value = '安信稳健聚申一年持有期混合A'
$newSheet.Range('D6').Value = '11.80'
$newSheet.Range('E6').Value = '40.12'
$newSheet.Range('F6').Value = '3.29'
$newSheet.Range('G6').Value = '0.3882'
$newSheet.Range('A6').Value = 4
$newSheet.Range('H6').Value = 5
$newSheet.Range('B7').Value = '008809'
$newSheet.Range('C7').Value = '安信民稳增长混合A'
$newSheet.Range('D7').Value = '9.19'
$newSheet.Range('E7').Value = '49.92'
$newSheet.Range('F7').Value = '4.12'
$newSheet.Range('G7').Value = '0.3786'
$newSheet.Range('A7').Value = 5
$newSheet.Range('H7').Value = 5
$newSheet.Range('B8').Value = '010865'
$newSheet.Range('C8').Value = '泓德卓远混合C'
$newSheet.Range('D8').Value = '11.56'
$newSheet.Range('E8').Value = '92.60'
$newSheet.Range('F8').Value = '2.83'
$newSheet.Range('G8').Value = '0.3271'
$newSheet.Range('A8').Value = 6
$newSheet.Range('H8').Value = 10
$newSheet.Range('B9').Value = '008810'
$newSheet.Range('C9').Value = '安信民稳增长混合C'
$newSheet.Range('D9').Value = '6.44'
$newSheet.Range('E9').Value = '49.92'
$newSheet.Range('F9').Value = '4.12'
$newSheet.Range('G9').Value = '0.2653'
$newSheet.Range('A9').Value = 7
$newSheet.Range('H9').Value = 5
$newSheet.Range('B10').Value = '015978'
$newSheet.Range('C10').Value = '安信恒鑫增强债券A'
$newSheet.Range('D10').Value = '16.87'
$newSheet.Range('E10').Value = '20.04'
$newSheet.Range('F10').Value = '1.43'
$newSheet.Range('G10').Value = '0.2412'
$newSheet.Range('A10').Value = 8
$newSheet.Range('H10').Value = 7
$newSheet.Range('B11').Value = '015979'
$newSheet.Range('C11').Value = '安信恒鑫增强债券C'
$newSheet.Range('D11').Value = '14.63'
$newSheet.Range('E11').Value = '20.04'
$newSheet.Range('F11').Value = '1.43'
$newSheet.Range('G11').Value = '0.2092'
$newSheet.Range('A11').Value = 9
$newSheet.Range('H11').Value = 7
$newSheet.Range('B12').Value = '012702'
$newSheet.Range('C12').Value = '安信民安回报一年持有混合C'
$newSheet.Range('D12').Value = '10.36'
$newSheet.Range('E12').Value = '25.63'
$newSheet.Range('F12').Value = '1.93'
$newSheet.Range('G12').Value = '0.1999'
$newSheet.Range('A12').Value = 10
$newSheet.Range('H12').Value = 7
$newSheet.Range('B13').Value = '012250'
$newSheet.Range('C13').Value = '安信平衡增利混合A'
$newSheet.Range('D13').Value = '3.30'
$newSheet.Range('E13').Value = '63.31'
$newSheet.Range('F13').Value = '4.71'
$newSheet.Range('G13').Value = '0.1554'
$newSheet.Range('A13').Value = 11
$newSheet.Range('H13').Value = 7
$newSheet.Range('B14').Value = '671010'
$newSheet.Range('C14').Value = '西部利得策略优选混合A'
$newSheet.Range('D14').Value = '1.88'
$newSheet.Range('E14').Value = '92.90'
$newSheet.Range('F14').Value = '5.56'
$newSheet.Range('G14').Value = '0.1045'
$newSheet.Range('A14').Value = 12
$newSheet.Range('H14').Value = 7
$newSheet.Range('B15').Value = '012251'
$newSheet.Range('C15').Value = '安信平衡增利混合C'
$newSheet.Range('D15').Value = '1.94'
$newSheet.Range('E15').Value = '63.31'
$newSheet.Range('F15').Value = '4.71'
$newSheet.Range('G15').Value = '0.0914'
$newSheet.Range('A15').Value = 13
$newSheet.Range('H15').Value = 7
$newSheet.Range('B16').Value = '012701'
$newSheet.Range('C16').Value = '安信民安回报一年持有混合A'
$newSheet.Range('D16').Value = '3.65'
$newSheet.Range('E16').Value = '25.63'
$newSheet.Range('F16').Value = '1.93'
$newSheet.Range('G16').Value = '0.0704'
$newSheet.Range('A16').Value = 14
$newSheet.Range('H16').Value = 7
$newSheet.Range('B17').Value = '010661'
$newSheet.Range('C17').Value = '安信稳健聚申一年持有期混合C'
$newSheet.Range('D17').Value = '1.70'
$newSheet.Range('E17').Value = '40.12'
$newSheet.Range('F17').Value = '3.29'
$newSheet.Range('G17').Value = '0.0559'
$newSheet.Range('A17').Value = 15
$newSheet.Range('H17').Value = 5
$newSheet.Range('B18').Value = '012257'
$newSheet.Range('C18').Value = '安信丰穗一年持有混合C'
$newSheet.Range('D18').Value = '2.45'
$newSheet.Range('E18').Value = '25.87'
$newSheet.Range('F18').Value = '1.81'
$newSheet.Range('G18').Value = '0.0443'
$newSheet.Range('A18').Value = 16
$newSheet.Range('H18').Value = 7
$newSheet.Range('B19').Value = '011969'
$newSheet.Range('C19').Value = '建信港股通精选混合A'
$newSheet.Range('D19').Value = '0.62'
$newSheet.Range('E19').Value = '87.38'
$newSheet.Range('F19').Value = '4.93'
$newSheet.Range('G19').Value = '0.0306'
$newSheet.Range('A19').Value = 17
$newSheet.Range('H19').Value = 7
$newSheet.Range('B20').Value = '011060'
$newSheet.Range('C20').Value = '西部利得策略优选混合C'
$newSheet.Range('D20').Value = '0.48'
$newSheet.Range('E20').Value = '92.90'
$newSheet.Range('F20').Value = '5.56'
$newSheet.Range('G20').Value = '0.0267'
$newSheet.Range('A20').Value = 18
$newSheet.Range('H20').Value = 7
$newSheet.Range('B21').Value = '011970'
$newSheet.Range('C21').Value = '建信港股通精选混合C'
$newSheet.Range('D21').Value = '0.24'
$newSheet.Range('E21').Value = '87.38'
$newSheet.Range('F21').Value = '4.93'
$newSheet.Range('G21').Value = '0.0118'
$newSheet.Range('A21').Value = 19
$newSheet.Range('H21').Value = 7

Write-Host "2022-Q4 sheet added and 总计 summary updated"
